$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.28521728515625
$ws.Range("B1").Value = 1.459593415260315
$ws.Range("C1").Value = 3.732053756713867
$ws.Range("D1").Value = 3.510660171508789
$ws.Range("E1").Value = 1.016089677810669
